# The deck ships two theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" colours (currently only wired to the Notes Master)
#   ppt/theme/theme2.xml -> "Integral" colours   (the live design used by the Slide Master)
#
# The authored change swaps the two themes' contents, so the deck's live design
# becomes "Office Theme" (and the Notes Master's theme becomes "Integral"). The
# font scheme and format scheme are identical between the two themes, so the only
# real content change is the 12 theme colours (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).
#
# Re-point the presentation's live colour scheme (Slide Master's Theme) from the
# "Integral" palette to the "Office Theme" palette to realise that swap.

function ToOleColor($r, $g, $b) {
    return ($b * 65536) + ($g * 256) + $r
}

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

$colors.Item(1).RGB  = ToOleColor 0x00 0x00 0x00   # dk1      000000
$colors.Item(2).RGB  = ToOleColor 0xFF 0xFF 0xFF   # lt1      FFFFFF
$colors.Item(3).RGB  = ToOleColor 0x44 0x54 0x6A   # dk2      44546A
$colors.Item(4).RGB  = ToOleColor 0xE7 0xE6 0xE6   # lt2      E7E6E6
$colors.Item(5).RGB  = ToOleColor 0x5B 0x9B 0xD5   # accent1  5B9BD5
$colors.Item(6).RGB  = ToOleColor 0xED 0x7D 0x31   # accent2  ED7D31
$colors.Item(7).RGB  = ToOleColor 0xA5 0xA5 0xA5   # accent3  A5A5A5
$colors.Item(8).RGB  = ToOleColor 0xFF 0xC0 0x00   # accent4  FFC000
$colors.Item(9).RGB  = ToOleColor 0x44 0x72 0xC4   # accent5  4472C4
$colors.Item(10).RGB = ToOleColor 0x70 0xAD 0x47   # accent6  70AD47
$colors.Item(11).RGB = ToOleColor 0x05 0x63 0xC1   # hlink    0563C1
$colors.Item(12).RGB = ToOleColor 0x95 0x4F 0x72   # folHlink 954F72
